# Reverse the order of the comma-separated "Recorded By" names in column G
# for every data row. Cells with a single value (no comma) are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*, *") {
        $parts = $val -split ", "
        $count = $parts.Count

        $reversedParts = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }

        $newVal = $reversedParts -join ", "
        $cell.Value2 = $newVal
    }
}
